$d = $word.ActiveDocument

# --- 1. Register the "Hyperlink" character style with the right formatting ---
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = "DefaultParagraphFont"
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Underline = 1
$hlStyle.Font.Color = 0xC16305   # BGR encoding of RGB 0563C1 (hyperlink theme color)

# --- 2. The original (only) paragraph carries the _GoBack bookmark. We leave
#        it completely untouched structurally and just insert 10 brand-new
#        blank paragraphs in front of it (1 plain-text + 9 hyperlinks), plus
#        1 new blank paragraph after it (plain text). The bookmark therefore
#        never needs to move. ---
$origPara = $d.Paragraphs.Item(1)
$origRange = $origPara.Range
$insertPoint = $origRange.Duplicate
$insertPoint.Collapse(1)   # start of the original paragraph

for ($i = 0; $i -lt 10; $i++) {
    $insertPoint.InsertParagraphBefore()
}

# The original paragraph (with the bookmark) is now #11; paragraphs 1-10 are
# the fresh blank ones we just inserted.

# --- helper: fill an already-blank paragraph (by 1-based index) with plain
#             text, replacing its empty run. ---
function Set-PlainParagraphText($index, $text) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.Text = $text
}

# --- helper: fill an already-blank paragraph (by 1-based index) with a
#             hyperlink whose display text is the url itself. ---
function Set-HyperlinkParagraphText($index, $url) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $r.Text = $url
    $p2 = $d.Paragraphs.Item($index)
    $r2 = $p2.Range
    $textRange = $d.Range($r2.Start, $r2.End - 1)
    $d.Hyperlinks.Add($textRange, $url, $null, $null, $null) | Out-Null
}

# Paragraph 1: plain stackoverflow link text
Set-PlainParagraphText 1 "http://stackoverflow.com/questions/28601663/how-to-retrieve-specific-node-from-firebase-database-in-android"

# Paragraphs 2-10: the nine hyperlinks
$urls = @(
  "https://dinosaur-facts.firebaseio.com/",
  "https://www.firebase.com/docs/web/api/query/equalto.html",
  "https://www.firebase.com/docs/android/guide/retrieving-data.html",
  "https://www.firebase.com/docs/web/guide/saving-data.html",
  "https://docs-examples.firebaseio.com/android/saving-data/fireblog",
  "https://www.firebase.com/docs/android/guide/saving-data.html",
  "https://www.firebase.com/docs/ios/guide/user-auth.html",
  "https://www.firebase.com/docs/android/guide/setup.html",
  "https://www.firebase.com/docs/web/guide/login/password.html"
)

$paraIndex = 2
foreach ($u in $urls) {
    Set-HyperlinkParagraphText $paraIndex $u
    $paraIndex = $paraIndex + 1
}

# --- 3. Paragraph 11 is the original bookmarked paragraph; update its text
#        in place (bookmark stays attached automatically). ---
$bcIndex = 11
$bcPara = $d.Paragraphs.Item($bcIndex)
$bcRange = $bcPara.Range
$bcTextOnly = $d.Range($bcRange.Start, $bcRange.End - 1)
$bcTextOnly.Text = "https://broadcast11.firebaseio.com/"

# --- 4. Append one new trailing plain-text paragraph after it. ---
$bcPara2 = $d.Paragraphs.Item($bcIndex)
$bcRange2 = $bcPara2.Range
$bcRange2.Collapse(0)
$bcRange2.InsertParagraphAfter()
$lastIndex = $bcIndex + 1
Set-PlainParagraphText $lastIndex "https://www.firebase.com/docs/android/guide/retrieving-data.html"

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ($i.ToString() + ": " + $d.Paragraphs.Item($i).Range.Text)
}
Write-Output ("Bookmark _GoBack exists=" + $d.Bookmarks.Exists("_GoBack"))
if ($d.Bookmarks.Exists("_GoBack")) {
    $gb = $d.Bookmarks.Item("_GoBack")
    Write-Output ("Bookmark range text=[" + $gb.Range.Text + "]")
}
